# "Actualizacion desde MV -datos-": append 6 new daily EMBI spread rows
# (25-10-2021 .. 01-11-2021) as rows 205-210 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date labels that must be stored as plain TEXT, exactly like
# every other "Serie" cell in the sheet above them. A direct
#   $ws.Cells.Item($r, 1).Value = "01-11-2021"
# is unsafe for day-of-month <= 12 (e.g. "01-11-2021"): Excel's smart text
# entry treats that as ambiguous and silently reinterprets it as a date
# serial (style + numeric value), which does not match the source data.
#
# Workaround: render the label as a text-formula result in a scratch cell
# (well outside the sheet's real A1:P210 data), Copy it, then
# PasteSpecial(values-only) into the destination. A values-only paste
# carries over the source's already-String-typed content verbatim (no
# text->date reparsing) and leaves the destination cell's style alone (stays
# on the sheet's default style, matching the unstyled cells in the diff).
$dates = @(
    "25-10-2021",
    "26-10-2021",
    "27-10-2021",
    "28-10-2021",
    "29-10-2021",
    "01-11-2021"
)

# Columns B..P for each new row, in order.
$rows = @(
    @(325,   397, 199, 323, 86.7, 119.1, 34.1, 150.6, 526.5, 1667.3, 322.5, 288, 156, 354, 159),
    @(325.2, 397, 202, 318, 86.7, 119.1, 37,   151.6, 510.9, 1664.2, 321.6, 288, 156, 354, 164),
    @(327.7, 398, 205, 321, 85.59999999999999, 118.2, 38.5, 157.3, 509.1, 1668.7, 320.9, 291, 155, 354, 166),
    @(324.7, 397, 200, 319, 85,   116.6, 34.1, 154.4, 509.4, 1693.6, 322,   294, 154, 348, 171),
    @(326.8, 401, 201, 322, 82.2, 117.1, 39.5, 156.5, 514.6, 1712.1, 337.8, 302, 161, 353, 181),
    @(326.4, 404, 198, 319, 81,   114,   35.6, 155.9, 507.2, 1728.4, 341.4, 303, 161, 357, 181)
)

$scratch = $ws.Range("Z1")
$startRow = 205

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i

    # --- Column A: date label, forced text via the scratch-cell trick. ---
    $scratch.Formula = "=""" + $dates[$i] + """"
    $scratch.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()

    # --- Columns B..P: plain numeric values. ---
    $values = $rows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $values[$c]
    }
}
